# Update the cryptocurrency price/volume data to the latest snapshot.
# Mirrors the per-cell content changes captured in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($sheet, [string]$addr, [string]$text)
    $cell = $sheet.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "1.00", "0.0912")
    # are kept verbatim instead of being reinterpreted as numbers, then
    # drop back to the default style so no stray formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" '29.669.62'
Set-CellText $ws "E2" '  +1.60%  '
Set-CellText $ws "D3" '1.602.72'
Set-CellText $ws "E3" '  +1.35%  '
Set-CellText $ws "D4" '1.00'
Set-CellText $ws "D5" '212.23'
Set-CellText $ws "E5" '  -0.44%  '
Set-CellText $ws "E7" '  +0.12%  '
Set-CellText $ws "D8" '27.83'
Set-CellText $ws "E8" '  +5.61%  '
Set-CellText $ws "E10" '  +1.41%  '
Set-CellText $ws "D11" '0.0912'
Set-CellText $ws "E11" '  +0.74%  '
Set-CellText $ws "D12" '1.831.13'
Set-CellText $ws "E12" '  +1.36%  '
Set-CellText $ws "D13" '1.612.16'
Set-CellText $ws "E13" '  +1.91%  '
Set-CellText $ws "D14" '0.545'
Set-CellText $ws "E14" '  +3.93%  '
Set-CellText $ws "D15" '29.660.14'
Set-CellText $ws "E15" '  +1.43%  '
Set-CellText $ws "E16" '  +0.85%  '
Set-CellText $ws "D17" '63.96'
Set-CellText $ws "E17" '  +1.81%  '
Set-CellText $ws "D18" '242.97'
Set-CellText $ws "E18" '  +2.02%  '
Set-CellText $ws "D19" '7.77'
Set-CellText $ws "E19" '  +4.18%  '
Set-CellText $ws "E20" '  +1.15%  '
Set-CellText $ws "E21" '  +0.13%  '
Set-CellText $ws "D23" '9.43'
Set-CellText $ws "E23" '  +2.67%  '
Set-CellText $ws "E24" '  -0.17%  '
Set-CellText $ws "D25" '155.37'
Set-CellText $ws "E25" '  +0.74%  '
Set-CellText $ws "D26" '15.46'
Set-CellText $ws "E26" '  +1.93%  '
Set-CellText $ws "E27" '  +0.25%  '
Set-CellText $ws "E28" '  +0.88%  '
Set-CellText $ws "E29" '  +0.13%  '
Set-CellText $ws "E30" '  +2.61%  '
Set-CellText $ws "E31" '  -0.33%  '
Set-CellText $ws "E32" '  +0.37%  '
Set-CellText $ws "D33" '3.20'
Set-CellText $ws "E33" '  +3.62%  '
Set-CellText $ws "D34" '1.425.85'
Set-CellText $ws "E34" '  -0.03%  '
Set-CellText $ws "E35" '  +3.39%  '
Set-CellText $ws "B36" 'TrustWalletToken'
Set-CellText $ws "C36" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText $ws "D36" '1.03'
Set-CellText $ws "E36" '  -0.72%  '
Set-CellText $ws "B37" 'MXToken'
Set-CellText $ws "C37" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws "D37" '2.89'
Set-CellText $ws "E37" '  +4.41%  '
Set-CellText $ws "E38" '  -0.24%  '
Set-CellText $ws "E39" '  +3.04%  '
Set-CellText $ws "D40" '58.39'
Set-CellText $ws "E40" '  +6.39%  '
Set-CellText $ws "E41" '  +2.33%  '
Set-CellText $ws "E42" '  +5.98%  '
Set-CellText $ws "E43" '  +0.28%  '
Set-CellText $ws "E44" '  +3.12%  '
Set-CellText $ws "E45" '  +0.10%  '
Set-CellText $ws "D46" '66.43'
Set-CellText $ws "E46" '  +2.89%  '
Set-CellText $ws "D47" '0.977'
Set-CellText $ws "E47" '  +16.48%  '
Set-CellText $ws "E48" '  +0.10%  '
Set-CellText $ws "D49" '1.742.12'
Set-CellText $ws "E49" '  +1.32%  '
Set-CellText $ws "D50" '86.78'
Set-CellText $ws "E50" '  +1.57%  '
Set-CellText $ws "D51" '0.0₆0105'
Set-CellText $ws "E51" '  +4.01%  '
